$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (coin names, links, prices, and 1h volume/change %)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.683.16'
$ws.Range("E2").Value = '  -1.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.378.06'
$ws.Range("E3").Value = '  -2.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.83'
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.17'
$ws.Range("E6").Value = '  -4.57%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.377.35'
$ws.Range("E8").Value = '  -2.04%  '

$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.48'
$ws.Range("E10").Value = '  -3.75%  '

$ws.Range("E11").Value = '  -0.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.955.79'
$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.124'
$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.03'
$ws.Range("E15").Value = '  +0.38%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.418.76'
$ws.Range("E16").Value = '  -0.95%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000170'
$ws.Range("E17").Value = '  -2.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.833.45'
$ws.Range("E18").Value = '  -1.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.26'
$ws.Range("E19").Value = '  -0.91%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.09'
$ws.Range("E20").Value = '  -1.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.04'
$ws.Range("E21").Value = '  -4.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.67'
$ws.Range("E22").Value = '  +0.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.559'
$ws.Range("E23").Value = '  -1.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.33'
$ws.Range("E24").Value = '  +1.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("E25").Value = '  -0.35%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000117'
$ws.Range("E26").Value = '  -4.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.522.45'
$ws.Range("E27").Value = '  -2.07%  '

$ws.Range("E28").Value = '  -0.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("E30").Value = '  -5.11%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.06'
$ws.Range("E31").Value = '  -2.21%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.16'
$ws.Range("E32").Value = '  -0.41%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.44'
$ws.Range("E33").Value = '  -3.75%  '

$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.70'
$ws.Range("E35").Value = '  -1.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.92'
$ws.Range("E36").Value = '  -2.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.408.07'
$ws.Range("E37").Value = '  -1.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '167.04'
$ws.Range("E38").Value = '  +0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.01'
$ws.Range("E39").Value = '  -3.86%  '

$ws.Range("E40").Value = '  -3.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0775'
$ws.Range("E41").Value = '  -1.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.90'
$ws.Range("E42").Value = '  +3.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.782'
$ws.Range("E43").Value = '  -1.77%  '

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.45'
$ws.Range("E45").Value = '  -0.98%  '

$ws.Range("E46").Value = '  -1.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.68'
$ws.Range("E47").Value = '  -2.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.534.70'
$ws.Range("E48").Value = '  -2.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.12'
$ws.Range("E49").Value = '  -4.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.81'
$ws.Range("E50").Value = '  -2.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.85'
$ws.Range("E51").Value = '  -2.02%  '
